$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.773.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "'3.510.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'614.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.30%  "
$ws.Range("D6").Value = "'191.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("D10").Value = "'0.664"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "'53.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "'4.053.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "'621.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.16%  "
$ws.Range("D16").Value = "'69.890.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'18.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'3.522.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "'0.991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'109.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.14%  "
$ws.Range("D23").Value = "'17.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").Value = "'3.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.06%  "
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'9.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("D29").Value = "'34.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.35%  "
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("D31").Value = "'12.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("D34").Value = "'63.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "'3.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "'3.676.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'522.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("D40").Value = "'0.394"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").Value = "'36.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("D42").Value = "'0.0₃0775"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "'0.0470"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "'132.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -4.29%  "
